$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (EMSO / national-ID numbers) needs to be stored as Text so the
# leading zeros / values survive verbatim - set the number format BEFORE
# writing the values.
$ws.Range("C2:C9").NumberFormat = "@"

# New EMSO values for rows 2-9
$ws.Range("C2").Value = "0203952500137"
$ws.Range("C3").Value = "1308959500124"
$ws.Range("C4").Value = "2203962505231"
$ws.Range("C5").Value = "1809955500218"
$ws.Range("C6").Value = "2710963500313"
$ws.Range("C7").Value = "3107964505276"
$ws.Range("C8").Value = "2811000500017"
$ws.Range("C9").Value = "1402001505453"

# Gender formula now reads the 10th character (instead of the 9th) and
# compares it against "5" (instead of testing the 9th char against "0"),
# with the "moski"/"zenski" branch results swapped.
$ws.Range("H3:H9").Formula = '=IF(MID(C3, 10, 1)>= "5", "ženski", "moški")'
$ws.Range("H2").Formula = '=IF(MID(C2, 10, 1)>= "5", "ženski", "moški")'

# Move the active selection to C2
$ws.Range("C2").Select()
